$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column A; remaining columns (B:F) shift left to (A:E)
$ws.Range("A1").EntireColumn.Delete()
